$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price records arrive and are inserted above the existing
# (older) Kiwi price rows, pushing the rest of the table down by two rows.
$ws.Rows("164:165").Insert()

# Row 164 - Terminal Hortofrutícola Agro Chillán, Kiwi Hayward, Primera
$ws.Range("A164").Value = 7
$ws.Range("B164").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C164").Value = "Ñuble"
$ws.Range("D164").Value = 44769
$ws.Range("E164").Value = 16
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100101
$ws.Range("H164").Value = "Berries"
$ws.Range("I164").Value = 100101007
$ws.Range("J164").Value = "Kiwi"
$ws.Range("K164").Value = "Hayward"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 120
$ws.Range("N164").Value = 6500
$ws.Range("O164").Value = 7000
$ws.Range("P164").Value = 6750
$ws.Range("Q164").Value = "`$/bandeja 18 kilos"
$ws.Range("R164").Value = "Provincia de Curicó"
$ws.Range("S164").Value = 375
$ws.Range("T164").Value = 18

# Row 165 - Terminal Hortofrutícola Agro Chillán, Kiwi Hayward, Segunda
$ws.Range("A165").Value = 7
$ws.Range("B165").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C165").Value = "Ñuble"
$ws.Range("D165").Value = 44769
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100101
$ws.Range("H165").Value = "Berries"
$ws.Range("I165").Value = 100101007
$ws.Range("J165").Value = "Kiwi"
$ws.Range("K165").Value = "Hayward"
$ws.Range("L165").Value = "Segunda"
$ws.Range("M165").Value = 120
$ws.Range("N165").Value = 5500
$ws.Range("O165").Value = 6000
$ws.Range("P165").Value = 5750
$ws.Range("Q165").Value = "`$/bandeja 18 kilos"
$ws.Range("R165").Value = "Provincia de Curicó"
$ws.Range("S165").Value = 319
$ws.Range("T165").Value = 18
